$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

# Row 1, Col 1: "59 x 99" -> "49 x 55"
$c = $t.Cell(1,1)
$c.Range.Text = "49 x 55" + $nl + "  5    5" + $nl + "  ----" + $nl + "4|    |" + $nl + "9|    |"

# Row 1, Col 2: "16 x 70" -> "95 x 42"
$c = $t.Cell(1,2)
$c.Range.Text = "95 x 42" + $nl + "  4    2" + $nl + "  ----" + $nl + "9|    |" + $nl + "5|    |"

# Row 1, Col 3: "82 x 39" -> "67 x 46"
$c = $t.Cell(1,3)
$c.Range.Text = "67 x 46" + $nl + "  4    6" + $nl + "  ----" + $nl + "6|    |" + $nl + "7|    |"

# Row 2, Col 1: "75 x 43" -> "74 x 90"
$c = $t.Cell(2,1)
$c.Range.Text = "74 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "7|    |" + $nl + "4|    |"

# Row 2, Col 2: "27 x 59" -> "59 x 27"
$c = $t.Cell(2,2)
$c.Range.Text = "59 x 27" + $nl + "  2    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"

# Row 2, Col 3: "59 x 61" -> "40 x 73"
$c = $t.Cell(2,3)
$c.Range.Text = "40 x 73" + $nl + "  7    3" + $nl + "  ----" + $nl + "4|    |" + $nl + "0|    |"

# Row 3, Col 1: "57 x 83" -> "99 x 14"
$c = $t.Cell(3,1)
$c.Range.Text = "99 x 14" + $nl + "  1    4" + $nl + "  ----" + $nl + "9|    |" + $nl + "9|    |"

# Row 3, Col 2: "29 x 97" -> "12 x 13"
$c = $t.Cell(3,2)
$c.Range.Text = "12 x 13" + $nl + "  1    3" + $nl + "  ----" + $nl + "1|    |" + $nl + "2|    |"

# Row 3, Col 3: "82 x 98" -> "29 x 51"
$c = $t.Cell(3,3)
$c.Range.Text = "29 x 51" + $nl + "  5    1" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"

# Row 4, Col 1: "29 x 60" -> "50 x 93"
$c = $t.Cell(4,1)
$c.Range.Text = "50 x 93" + $nl + "  9    3" + $nl + "  ----" + $nl + "5|    |" + $nl + "0|    |"

# Row 4, Col 2: "71 x 25" -> "29 x 19"
$c = $t.Cell(4,2)
$c.Range.Text = "29 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"

# Row 4, Col 3: "84 x 51" -> "51 x 67"
$c = $t.Cell(4,3)
$c.Range.Text = "51 x 67" + $nl + "  6    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "1|    |"

# Row 5, Col 1: "91 x 36" -> "66 x 62"
$c = $t.Cell(5,1)
$c.Range.Text = "66 x 62" + $nl + "  6    2" + $nl + "  ----" + $nl + "6|    |" + $nl + "6|    |"

# Row 5, Col 2: "63 x 93" -> "96 x 67"
$c = $t.Cell(5,2)
$c.Range.Text = "96 x 67" + $nl + "  6    7" + $nl + "  ----" + $nl + "9|    |" + $nl + "6|    |"

# Row 5, Col 3: "13 x 88" -> "22 x 48"
$c = $t.Cell(5,3)
$c.Range.Text = "22 x 48" + $nl + "  4    8" + $nl + "  ----" + $nl + "2|    |" + $nl + "2|    |"
